# "DANH SÁCH NỢ" sheet: two new debtors added to the tracking table.
# The table previously had two spare blank rows (37-38) right before the
# "Điều khoản dịch vụ" / totals block (rows 39-43). We fill those two rows
# with the new debtors and insert 6 more blank rows to restore the spare
# capacity, then refresh the summary formulas below to cover the bigger
# range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DANH SÁCH NỢ")

# Make room: push the "Điều khoản dịch vụ" / totals block down by 6 rows,
# keeping rows 37-38 (currently blank) as the insertion point for the new
# debtors.
$ws.Rows("39:44").Insert()

# --- Row 37: Lê Ngọc Như Ý / Nạp quân huy ---------------------------------
$ws.Range("B37").Value = "Lê Ngọc Như Ý"
$ws.Range("C37").Value = "Nạp quân huy"
$ws.Range("D37").Value = 100000
$ws.Range("E37").Value = 0
$ws.Range("F37").Formula = "=(D37+I37)-E37"
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Formula = "=D37*H37"
$ws.Range("J37").Value = 46075
$ws.Range("K37").Value = 46083
$ws.Range("M37").Value = "Chưa trả đủ"

# --- Row 38: Liêu Thuận Phát / Nạp Free Fire ------------------------------
$ws.Range("B38").Value = "Liêu Thuận Phát"
$ws.Range("C38").Value = "Nạp Free Fire"
$ws.Range("D38").Value = 505000
$ws.Range("E38").Value = 0
$ws.Range("F38").Formula = "=(D38+I38)-E38"
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("I38").Formula = "=D38*H38"
$ws.Range("J38").Value = 46075
$ws.Range("K38").Value = 46083
$ws.Range("M38").Value = "Chưa trả đủ"

# Match the style used by the existing data rows (copy the row-36 look
# down onto the two freshly-filled rows).
$ws.Range("A36:M36").Copy()
$ws.Range("A37:M38").PasteSpecial(-4122)
$ws.Range("B37").Value = "Lê Ngọc Như Ý"
$ws.Range("C37").Value = "Nạp quân huy"
$ws.Range("D37").Value = 100000
$ws.Range("E37").Value = 0
$ws.Range("F37").Formula = "=(D37+I37)-E37"
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Formula = "=D37*H37"
$ws.Range("J37").Value = 46075
$ws.Range("K37").Value = 46083
$ws.Range("M37").Value = "Chưa trả đủ"
$ws.Range("B38").Value = "Liêu Thuận Phát"
$ws.Range("C38").Value = "Nạp Free Fire"
$ws.Range("D38").Value = 505000
$ws.Range("E38").Value = 0
$ws.Range("F38").Formula = "=(D38+I38)-E38"
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("I38").Formula = "=D38*H38"
$ws.Range("J38").Value = 46075
$ws.Range("K38").Value = 46083
$ws.Range("M38").Value = "Chưa trả đủ"

# --- Refresh the totals block (now at rows 45-52) -------------------------
$ws.Range("F46").Formula = "=SUM(D2:D44)"
$ws.Range("F47").Formula = "=SUM(E2:E44)"
$ws.Range("F49").Formula = "=SUM(F2:F44)"
$ws.Range("E51").Formula = "=1000000+1000000+1000000+3000000-1500000+300000+600000"
$ws.Range("E52").Formula = "=135000+300000"

# --- Move the "Điều khoản dịch vụ" hyperlink along with its cell ---------
$ws.Range("D39").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D45"), "https://tinyurl.com/dieukhoan02112026")

# --- Keep the autofilter / selection in step with the bigger table -------
$ws.AutoFilterMode = $false
$ws.Range("A1:M48").AutoFilter()
$ws.Range("H45").Select()
$excel.ActiveWindow.ScrollRow = 25

# --- "THONG KE NAP " sheet: log the same two payments in its ledger ------
$ws2 = $wb.Worksheets.Item("THONG KE NAP ")
$ws2.Range("A291").Value = 46075
$ws2.Range("B291").Value = "Lê Ngọc Như Ý"
$ws2.Range("C291").Value = 100000
$ws2.Range("D291").Value = "Nạp quân huy"

$ws2.Range("A292").Value = 46075
$ws2.Range("B292").Value = "Liêu Thuận Phát"
$ws2.Range("C292").Value = 505000
$ws2.Range("D292").Value = "Nạp Free Fire"

$ws2.Range("D293").Select()
